# Milwaukee_B team-specific matrix refresh: probabilities recomputed from
# games pulled March 7. Update the affected matrix cells (rows 2-19,
# columns B-S) on Sheet1 with the new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1794871794871795
$ws.Range("C2").Value = 0.6011396011396012
$ws.Range("J2").Value = 0.01994301994301994
$ws.Range("P2").Value = 0.1225071225071225
$ws.Range("S2").Value = 0.07692307692307693
$ws.Range("B3").Value = 0.0186046511627907
$ws.Range("C3").Value = 0.02325581395348837
$ws.Range("J3").Value = 0.04651162790697674
$ws.Range("P3").Value = 0.772093023255814
$ws.Range("S3").Value = 0.1395348837209302
$ws.Range("P4").Value = 0.78
$ws.Range("S4").Value = 0.22
$ws.Range("B6").Value = 0.07024793388429752
$ws.Range("D6").Value = 0.008264462809917356
$ws.Range("F6").Value = 0.04958677685950413
$ws.Range("J6").Value = 0.256198347107438
$ws.Range("O6").Value = 0.0371900826446281
$ws.Range("Q6").Value = 0.1694214876033058
$ws.Range("R6").Value = 0.09917355371900827
$ws.Range("S6").Value = 0.3099173553719008
$ws.Range("B7").Value = 0.1106382978723404
$ws.Range("D7").Value = 0.01702127659574468
$ws.Range("F7").Value = 0.06808510638297872
$ws.Range("J7").Value = 0.1063829787234043
$ws.Range("O7").Value = 0.02127659574468085
$ws.Range("Q7").Value = 0.1787234042553192
$ws.Range("R7").Value = 0.05531914893617021
$ws.Range("S7").Value = 0.4425531914893617
$ws.Range("B8").Value = 0.1333333333333333
$ws.Range("D8").Value = 0.02156862745098039
$ws.Range("F8").Value = 0.05098039215686274
$ws.Range("J8").Value = 0.1294117647058824
$ws.Range("O8").Value = 0.01176470588235294
$ws.Range("Q8").Value = 0.1725490196078431
$ws.Range("R8").Value = 0.103921568627451
$ws.Range("S8").Value = 0.3764705882352941
$ws.Range("B9").Value = 0.09677419354838709
$ws.Range("D9").Value = 0.01612903225806452
$ws.Range("F9").Value = 0.07258064516129033
$ws.Range("J9").Value = 0.1048387096774194
$ws.Range("O9").Value = 0.01209677419354839
$ws.Range("Q9").Value = 0.1975806451612903
$ws.Range("R9").Value = 0.09677419354838709
$ws.Range("S9").Value = 0.4032258064516129
$ws.Range("B10").Value = 0.09817671809256662
$ws.Range("D10").Value = 0.02103786816269285
$ws.Range("E10").Value = 0.001402524544179523
$ws.Range("F10").Value = 0.06591865357643759
$ws.Range("J10").Value = 0.1416549789621318
$ws.Range("O10").Value = 0.01122019635343618
$ws.Range("Q10").Value = 0.2047685834502104
$ws.Range("R10").Value = 0.09046283309957924
$ws.Range("S10").Value = 0.3653576437587658
$ws.Range("G11").Value = 0.1327683615819209
$ws.Range("J11").Value = 0.07344632768361582
$ws.Range("K11").Value = 0.172316384180791
$ws.Range("L11").Value = 0.5847457627118644
$ws.Range("S11").Value = 0.03672316384180791
$ws.Range("G12").Value = 0.6995305164319249
$ws.Range("J12").Value = 0.215962441314554
$ws.Range("K12").Value = 0.009389671361502348
$ws.Range("L12").Value = 0.02347417840375587
$ws.Range("S12").Value = 0.05164319248826291
$ws.Range("G13").Value = 0.7288135593220338
$ws.Range("J13").Value = 0.2033898305084746
$ws.Range("S13").Value = 0.06779661016949153
$ws.Range("F15").Value = 0.025
$ws.Range("H15").Value = 0.1541666666666667
$ws.Range("I15").Value = 0.09166666666666666
$ws.Range("J15").Value = 0.3625
$ws.Range("K15").Value = 0.05416666666666667
$ws.Range("M15").Value = 0.01666666666666667
$ws.Range("O15").Value = 0.05833333333333333
$ws.Range("S15").Value = 0.2375
$ws.Range("F16").Value = 0.02066115702479339
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3429752066115703
$ws.Range("K16").Value = 0.09504132231404959
$ws.Range("M16").Value = 0.02892561983471074
$ws.Range("O16").Value = 0.05785123966942149
$ws.Range("S16").Value = 0.1818181818181818
$ws.Range("F17").Value = 0.0198019801980198
$ws.Range("H17").Value = 0.1861386138613861
$ws.Range("I17").Value = 0.1069306930693069
$ws.Range("J17").Value = 0.3920792079207921
$ws.Range("K17").Value = 0.07920792079207921
$ws.Range("M17").Value = 0.0198019801980198
$ws.Range("N17").Value = 0.00198019801980198
$ws.Range("O17").Value = 0.04158415841584159
$ws.Range("S17").Value = 0.1524752475247525
$ws.Range("F18").Value = 0.0205761316872428
$ws.Range("H18").Value = 0.168724279835391
$ws.Range("I18").Value = 0.1069958847736626
$ws.Range("J18").Value = 0.3950617283950617
$ws.Range("K18").Value = 0.1069958847736626
$ws.Range("M18").Value = 0.01234567901234568
$ws.Range("O18").Value = 0.06584362139917696
$ws.Range("S18").Value = 0.1234567901234568
$ws.Range("F19").Value = 0.01345895020188425
$ws.Range("H19").Value = 0.2032301480484522
$ws.Range("I19").Value = 0.08209959623149395
$ws.Range("J19").Value = 0.3391655450874832
$ws.Range("K19").Value = 0.1238223418573351
$ws.Range("M19").Value = 0.02422611036339166
$ws.Range("N19").Value = 0.0006729475100942127
$ws.Range("O19").Value = 0.07200538358008075
$ws.Range("S19").Value = 0.1413189771197847
